# Updates cryptos list data (Coin/Link/Price/Volume) per the Jul 22 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin (B), Link (C), Price (D), Volume1h (E), forceTextPrice (whether D must be kept as text)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "29.892.15", "  +0.20%  ", 0),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.889.01", "  -0.07%  ", 0),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.000", "  +0.04%  ", 1),
    @(5, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.7643", "  -1.42%  ", 1),
    @(6, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "242.77", "  -0.54%  ", 1),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  -0.05%  ", 1),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3130", "  -0.15%  ", 1),
    @(9, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "25.68", "  +1.57%  ", 1),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07144", "  -2.74%  ", 1),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.08565", "  +5.18%  ", 1),
    @(12, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.7629", "  -0.29%  ", 1),
    @(13, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.371", "  -1.57%  ", 1),
    @(14, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.856.85", "  -0.58%  ", 0),
    @(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "93.86", "  +0.85%  ", 1),
    @(16, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.182", "  -0.10%  ", 1),
    @(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "29.643.07", "  -0.61%  ", 0),
    @(18, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.78", "  -0.98%  ", 1),
    @(19, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "244.00", "  -0.62%  ", 1),
    @(20, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007802", "  -0.55%  ", 1),
    @(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9991", "  -0.09%  ", 1),
    @(22, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.995", "  -1.92%  ", 1),
    @(23, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  +0.09%  ", 1),
    @(24, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1617", "  +2.82%  ", 1),
    @(25, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.412", "  -0.07%  ", 1),
    @(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "162.77", "  +0.85%  ", 1),
    @(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.77", "  +0.07%  ", 1),
    @(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.034", "  -0.20%  ", 1),
    @(29, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.522", "  +4.90%  ", 1),
    @(30, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.539", "  -0.27%  ", 1),
    @(31, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.486", "  +0.30%  ", 1),
    @(32, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.115", "  +0.94%  ", 1),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05436", "  -2.32%  ", 1),
    @(34, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.240", "  -0.43%  ", 1),
    @(35, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7431", "  -1.65%  ", 1),
    @(36, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9988", "  +0.16%  ", 1),
    @(37, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.697", "  +2.25%  ", 1),
    @(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01949", "  +0.96%  ", 1),
    @(39, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.781", "  +0.26%  ", 1),
    @(40, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4471", "  +0.55%  ", 1),
    @(41, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.104.71", "  -3.57%  ", 0),
    @(42, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.082", "  +2.23%  ", 1),
    @(43, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "73.05", "  -0.51%  ", 1),
    @(44, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8535", "  +0.25%  ", 1),
    @(45, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.000", "  +0.05%  ", 1),
    @(46, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "102.82", "  +1.01%  ", 1),
    @(47, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.869", "  -1.52%  ", 1),
    @(48, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.656", "  +2.17%  ", 1),
    @(49, "SynthetixNetwork", "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx", "3.025", "  -3.06%  ", 1),
    @(50, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06087", "  +0.73%  ", 1),
    @(51, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.992.10", "  -1.62%  ", 0)
)

foreach ($entry in $data) {
    $rowNum      = $entry[0]
    $coin        = $entry[1]
    $link        = $entry[2]
    $price       = $entry[3]
    $volume      = $entry[4]
    $forceText   = $entry[5]

    $bCell = $ws.Cells.Item($rowNum, 2)
    $cCell = $ws.Cells.Item($rowNum, 3)
    $dCell = $ws.Cells.Item($rowNum, 4)
    $eCell = $ws.Cells.Item($rowNum, 5)

    $bCell.Value2 = $coin
    $cCell.Value2 = $link

    if ($forceText -eq 1) {
        # Keep values such as "244.00" or "0.000007802" stored as text, matching the source export,
        # instead of letting Excel coerce them into numbers.
        $dCell.NumberFormat = "@"
    }
    $dCell.Value2 = $price
    $eCell.Value2 = $volume
}

Write-Host "Done updating crypto rows."